# Insert a new weekly price record as row 85 ("Hortaliza, Femacal de La
# Calera - Berenjena"), pushing the existing rows 85-138 down to 86-139.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 85..138 down to 86..139, leaving a blank row 85 to populate.
$ws.Rows("85:85").Insert()

# Populate the newly inserted row 85 with the new record's data.
$ws.Cells.Item(85, 1).Value = 3
$ws.Cells.Item(85, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44452
$ws.Cells.Item(85, 5).Value = 5
$ws.Cells.Item(85, 6).Value = 100112001
$ws.Cells.Item(85, 7).Value = "Berenjena"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 73
$ws.Cells.Item(85, 11).Value = 9500
$ws.Cells.Item(85, 12).Value = 10000
$ws.Cells.Item(85, 13).Value = 9760
$ws.Cells.Item(85, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 163
$ws.Cells.Item(85, 17).Value = 60
$ws.Cells.Item(85, 18).Value = "Hortaliza"
